# "Generate Report for Handback"
#
# The CI report workbook is re-run after a handback completes: the
# "Status" column flips from "Ready for handoff" to
# "Handed back: in sync with en-US", the per-language sheets grow a
# "Latest Target File" / "Latest Handback File" pair of columns (with
# their own hyperlinks) now that a handback file exists, and the
# "Latest Handback DateTime" column is stamped with the real time the
# handback finished (instead of the zero-date placeholder).

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item(1)
$zhcn     = $wb.Worksheets.Item(2)
$dede     = $wb.Worksheets.Item(3)

$statusDone = "Handed back: in sync with en-US"

# --- Overview sheet: both language status columns flip to "handed back" ---
$overview.Range("B2").Value = $statusDone
$overview.Range("C2").Value = $statusDone
$overview.Range("B3").Value = $statusDone
$overview.Range("C3").Value = $statusDone

# --- Per-language sheets: Status column flips too ---
$zhcn.Range("C2").Value = $statusDone
$zhcn.Range("C3").Value = $statusDone
$dede.Range("C2").Value = $statusDone
$dede.Range("C3").Value = $statusDone

# --- Latest Handback DateTime (column H) gets a real timestamp ---
$zhcn.Range("H2").Value = "2016-03-20 14:32:08"
$zhcn.Range("H3").Value = "2016-03-20 14:32:08"
$dede.Range("H2").Value = "2016-03-20 14:32:15"
$dede.Range("H3").Value = "2016-03-20 14:32:15"

# --- New "Latest Target File" (F) / "Latest Handback File" (G) columns ---
# zh-cn
$zhcn.Range("F2").Value = "a.md"
$zhcn.Hyperlinks.Add($zhcn.Range("F2"), "https://github.com/OpenLocalizationTest/oltest/blob/f600f3eb0eda8033b2da1f02185e86826d4323d1/e2e/a.md", "", "", "a.md") | Out-Null
$zhcn.Range("F2").Style = "HyperLink"

$zhcn.Range("G2").Value = "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.zh-cn.xlf"
$zhcn.Hyperlinks.Add($zhcn.Range("G2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/0ec4006ef4431ffe8e9884457360302ba3b92f16/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.zh-cn.xlf", "", "", "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.zh-cn.xlf") | Out-Null
$zhcn.Range("G2").Style = "HyperLink"

$zhcn.Range("F3").Value = "a.md"
$zhcn.Hyperlinks.Add($zhcn.Range("F3"), "https://github.com/OpenLocalizationTest/oltest/blob/f600f3eb0eda8033b2da1f02185e86826d4323d1/e2e/a.md", "", "", "a.md") | Out-Null
$zhcn.Range("F3").Style = "HyperLink"

$zhcn.Range("G3").Value = "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.zh-cn.xlf"
$zhcn.Hyperlinks.Add($zhcn.Range("G3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/0ec4006ef4431ffe8e9884457360302ba3b92f16/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.zh-cn.xlf", "", "", "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.zh-cn.xlf") | Out-Null
$zhcn.Range("G3").Style = "HyperLink"

# de-de
$dede.Range("F2").Value = "a.md"
$dede.Hyperlinks.Add($dede.Range("F2"), "https://github.com/OpenLocalizationTest/oltest/blob/f600f3eb0eda8033b2da1f02185e86826d4323d1/e2e/a.md", "", "", "a.md") | Out-Null
$dede.Range("F2").Style = "HyperLink"

$dede.Range("G2").Value = "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.de-de.xlf"
$dede.Hyperlinks.Add($dede.Range("G2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/4157e2b3e32784db95b99f1cd95516e07c4fdefd/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.de-de.xlf", "", "", "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.de-de.xlf") | Out-Null
$dede.Range("G2").Style = "HyperLink"

$dede.Range("F3").Value = "a.md"
$dede.Hyperlinks.Add($dede.Range("F3"), "https://github.com/OpenLocalizationTest/oltest/blob/f600f3eb0eda8033b2da1f02185e86826d4323d1/e2e/a.md", "", "", "a.md") | Out-Null
$dede.Range("F3").Style = "HyperLink"

$dede.Range("G3").Value = "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.de-de.xlf"
$dede.Hyperlinks.Add($dede.Range("G3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/4157e2b3e32784db95b99f1cd95516e07c4fdefd/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.de-de.xlf", "", "", "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.de-de.xlf") | Out-Null
$dede.Range("G3").Style = "HyperLink"
